$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking strings
# like "1.00" or "0.999" are not silently coerced to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "98.244.80"
$ws.Cells.Item(2, 5).Value = "  +0.95%  "
$ws.Cells.Item(3, 4).Value = "3.327.19"
$ws.Cells.Item(3, 5).Value = "  -0.54%  "
$ws.Cells.Item(4, 5).Value = "  +0.13%  "
$ws.Cells.Item(5, 4).Value = "256.32"
$ws.Cells.Item(5, 5).Value = "  +3.47%  "
$ws.Cells.Item(6, 4).Value = "622.66"
$ws.Cells.Item(6, 5).Value = "  -0.20%  "
$ws.Cells.Item(7, 4).Value = "1.44"
$ws.Cells.Item(7, 5).Value = "  +27.07%  "
$ws.Cells.Item(8, 4).Value = "0.402"
$ws.Cells.Item(8, 5).Value = "  +3.38%  "
$ws.Cells.Item(9, 4).Value = "1.00"
$ws.Cells.Item(9, 5).Value = "  -0.01%  "
$ws.Cells.Item(10, 4).Value = "0.899"
$ws.Cells.Item(10, 5).Value = "  +12.71%  "
$ws.Cells.Item(11, 4).Value = "3.320.22"
$ws.Cells.Item(11, 5).Value = "  -0.84%  "
$ws.Cells.Item(12, 4).Value = "0.199"
$ws.Cells.Item(12, 5).Value = "  -0.51%  "
$ws.Cells.Item(13, 4).Value = "38.06"
$ws.Cells.Item(13, 5).Value = "  +6.35%  "
$ws.Cells.Item(14, 4).Value = "97.943.34"
$ws.Cells.Item(14, 5).Value = "  +0.27%  "
$ws.Cells.Item(15, 4).Value = "0.0000250"
$ws.Cells.Item(15, 5).Value = "  +0.82%  "
$ws.Cells.Item(16, 4).Value = "3.964.23"
$ws.Cells.Item(16, 5).Value = "  +2.88%  "
$ws.Cells.Item(17, 4).Value = "5.51"
$ws.Cells.Item(17, 5).Value = "  -0.38%  "
$ws.Cells.Item(18, 4).Value = "3.325.15"
$ws.Cells.Item(18, 5).Value = "  +0.79%  "
$ws.Cells.Item(19, 4).Value = "3.55"
$ws.Cells.Item(19, 5).Value = "  -2.98%  "
$ws.Cells.Item(20, 4).Value = "15.20"
$ws.Cells.Item(20, 5).Value = "  -0.43%  "
$ws.Cells.Item(21, 4).Value = "481.21"
$ws.Cells.Item(21, 5).Value = "  -2.64%  "
$ws.Cells.Item(22, 4).Value = "6.09"
$ws.Cells.Item(22, 5).Value = "  +2.44%  "
$ws.Cells.Item(23, 4).Value = "0.0000205"
$ws.Cells.Item(23, 5).Value = "  -3.18%  "
$ws.Cells.Item(24, 4).Value = "9.38"
$ws.Cells.Item(24, 5).Value = "  +0.85%  "
$ws.Cells.Item(25, 4).Value = "5.58"
$ws.Cells.Item(25, 5).Value = "  -2.12%  "
$ws.Cells.Item(26, 4).Value = "89.08"
$ws.Cells.Item(26, 5).Value = "  +0.34%  "
$ws.Cells.Item(27, 4).Value = "11.89"
$ws.Cells.Item(27, 5).Value = "  -2.55%  "
$ws.Cells.Item(28, 2).Value = "WrappedeETH"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(28, 4).Value = "3.509.81"
$ws.Cells.Item(28, 5).Value = "  +0.37%  "
$ws.Cells.Item(29, 2).Value = "Stellar"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(29, 4).Value = "0.291"
$ws.Cells.Item(29, 5).Value = "  +19.67%  "
$ws.Cells.Item(30, 5).Value = "  +0.05%  "
$ws.Cells.Item(31, 4).Value = "0.188"
$ws.Cells.Item(31, 5).Value = "  +3.49%  "
$ws.Cells.Item(32, 4).Value = "0.132"
$ws.Cells.Item(32, 5).Value = "  +6.02%  "
$ws.Cells.Item(33, 4).Value = "9.80"
$ws.Cells.Item(33, 5).Value = "  +4.59%  "
$ws.Cells.Item(34, 4).Value = "0.999"
$ws.Cells.Item(35, 4).Value = "27.80"
$ws.Cells.Item(35, 5).Value = "  +0.21%  "
$ws.Cells.Item(36, 2).Value = "Kaspa"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(36, 4).Value = "0.148"
$ws.Cells.Item(36, 5).Value = "  -5.63%  "
$ws.Cells.Item(37, 2).Value = "RenderToken"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(37, 4).Value = "7.17"
$ws.Cells.Item(37, 5).Value = "  -5.23%  "
$ws.Cells.Item(38, 4).Value = "1.93"
$ws.Cells.Item(38, 5).Value = "  -0.99%  "
$ws.Cells.Item(39, 2).Value = "Bittensor"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(39, 4).Value = "496.93"
$ws.Cells.Item(39, 5).Value = "  -2.85%  "
$ws.Cells.Item(40, 2).Value = "PolygonEcosystemToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Cells.Item(40, 4).Value = "0.460"
$ws.Cells.Item(40, 5).Value = "  +1.59%  "
$ws.Cells.Item(41, 4).Value = "24.88"
$ws.Cells.Item(41, 5).Value = "  +0.43%  "
$ws.Cells.Item(42, 2).Value = "MantraDAO"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Cells.Item(42, 4).Value = "3.68"
$ws.Cells.Item(42, 5).Value = "  +4.59%  "
$ws.Cells.Item(43, 2).Value = "Fetch.AI"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(43, 4).Value = "1.25"
$ws.Cells.Item(43, 5).Value = "  -3.07%  "
$ws.Cells.Item(44, 2).Value = "dogwifhat"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(44, 4).Value = "3.27"
$ws.Cells.Item(44, 5).Value = "  -1.00%  "
$ws.Cells.Item(45, 2).Value = "ARBITRUM"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(45, 4).Value = "0.794"
$ws.Cells.Item(45, 5).Value = "  +1.98%  "
$ws.Cells.Item(47, 2).Value = "Monero"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(47, 4).Value = "160.42"
$ws.Cells.Item(47, 5).Value = "  -0.86%  "
$ws.Cells.Item(48, 2).Value = "Stacks"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(48, 4).Value = "1.92"
$ws.Cells.Item(48, 5).Value = "  -2.05%  "
$ws.Cells.Item(49, 4).Value = "0.841"
$ws.Cells.Item(49, 5).Value = "  +5.42%  "
$ws.Cells.Item(50, 4).Value = "4.65"
$ws.Cells.Item(50, 5).Value = "  +1.07%  "
$ws.Cells.Item(51, 2).Value = "OKB"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(51, 4).Value = "45.80"
$ws.Cells.Item(51, 5).Value = "  +1.63%  "

# Restore original (General) formatting/style on column D now that
# the text values are committed, so no stray number format lingers.
$dRange.ClearFormats()
